# Insert a new weekly price record at row 243 of the data table.
# This shifts the existing rows 243:290 down to 244:291 (Excel keeps the
# D-column date style on the shifted cells automatically) and grows the
# sheet's used range to A1:R291.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("243:243").Insert()

$ws.Range("A243").Value = 7
$ws.Range("B243").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C243").Value = "Ñuble"
$ws.Range("D243").Value = 44951
$ws.Range("E243").Value = 16
$ws.Range("F243").Value = 100112043
$ws.Range("G243").Value = "Pepino ensalada"
$ws.Range("H243").Value = "Sin especificar"
$ws.Range("I243").Value = "Primera"
$ws.Range("J243").Value = 50
$ws.Range("K243").Value = 9000
$ws.Range("L243").Value = 9000
$ws.Range("M243").Value = 9000
$ws.Range("N243").Value = "`$/caja 80 unidades"
$ws.Range("O243").Value = "Región del Maule"
$ws.Range("P243").Value = 112
$ws.Range("Q243").Value = 80
$ws.Range("R243").Value = "Hortaliza"
